$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CasesTab): update the "query" column (B2) to the new case-detail query.
# Column C2 (StatQuery) becomes the shared Trials/Cases/Files count query.
$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
   WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
  WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# Write order matters for shared-string index assignment, so introduce each
# brand-new unique string in the same sequence the original authoring tool did:
# FilesTab, then the shared StatQuery text, then the Files query, then the Cases query.
$ws.Range("A3").Value = "FilesTab"
$ws.Range("C2").Value = $statQuery
$ws.Range("B3").Value = $filesQuery
$ws.Range("B2").Value = $casesQuery

$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = "TC01_Trials_Filter_Race-AmerIndAlask_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Trials_Filter_Race-AmerIndAlask_WebData.xlsx"

$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Row heights grow to fit the taller wrapped query text (row 3's query is long
# enough to hit Excel's row-height cap of 409.5 points).
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

$ws.Range("A2:C3").Select()
